$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (for new quarters: Dec-2018, Sep-2018)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy the number/date formatting from the (now shifted) old column D, now at F, into new D:E
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns with the new quarterly data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 478200
$ws.Range("E8").Value = 471600
$ws.Range("D9").Value = 357600
$ws.Range("E9").Value = 358500
$ws.Range("D10").Value = 120600
$ws.Range("E10").Value = 113100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 8100
$ws.Range("E15").Value = 8200
$ws.Range("D17").Value = 414500
$ws.Range("E17").Value = 416100
$ws.Range("D18").Value = 63700
$ws.Range("E18").Value = 55500
$ws.Range("D20").Value = 600
$ws.Range("E20").Value = 200
$ws.Range("D21").Value = 90000
$ws.Range("E21").Value = 80900
$ws.Range("D22").Value = 12600
$ws.Range("E22").Value = 12000
$ws.Range("D23").Value = 51700
$ws.Range("E23").Value = 43700
$ws.Range("D24").Value = 5000
$ws.Range("E24").Value = 10100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 46700
$ws.Range("E26").Value = 33600
$ws.Range("D27").Value = 46700
$ws.Range("E27").Value = 33600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -600
$ws.Range("E32").Value = -200
$ws.Range("D33").Value = 46700
$ws.Range("E33").Value = 33600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 46700
$ws.Range("E35").Value = 33600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 15500
$ws.Range("E41").Value = 21300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 131200
$ws.Range("E43").Value = 105000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 47300
$ws.Range("E45").Value = 48200
$ws.Range("D46").Value = 193900
$ws.Range("E46").Value = 174500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 597100
$ws.Range("E48").Value = 586400
$ws.Range("D49").Value = 1670600
$ws.Range("E49").Value = 1668800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 62600
$ws.Range("E52").Value = 67800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2524300
$ws.Range("E54").Value = 2497400
$ws.Range("D57").Value = 154200
$ws.Range("E57").Value = 160000
$ws.Range("D58").Value = 129000
$ws.Range("E58").Value = 104800
$ws.Range("D59").Value = 200600
$ws.Range("E59").Value = 187700
$ws.Range("D60").Value = 483800
$ws.Range("E60").Value = 452500
$ws.Range("D61").Value = 1036900
$ws.Range("E61").Value = 1039100
$ws.Range("D62").Value = 224200
$ws.Range("E62").Value = 225600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1744800
$ws.Range("E66").Value = 1717200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 193100
$ws.Range("E72").Value = 146500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 779500
$ws.Range("E76").Value = 780200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 46700
$ws.Range("E81").Value = 33600
$ws.Range("D83").Value = 25800
$ws.Range("E83").Value = 25200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 55100
$ws.Range("E89").Value = 51000
$ws.Range("D91").Value = -29900
$ws.Range("E91").Value = -23400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -44300
$ws.Range("E94").Value = -24500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -15900
$ws.Range("E100").Value = -18500
$ws.Range("D101").Value = 200
$ws.Range("E101").Value = -800
$ws.Range("D102").Value = -5000
$ws.Range("E102").Value = 7200

# Correct a handful of historical figures that were also revised in this update
$ws.Range("G20").Value = -400
$ws.Range("G21").Value = 79600
$ws.Range("G22").Value = 11100
$ws.Range("G32").Value = 400
$ws.Range("H89").Value = 47000
$ws.Range("H102").Value = -7100
